$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 3721112.21
$ws.Range("C7").Value = -16.24922836153059
$ws.Range("D7").Value = 3294
$ws.Range("E7").Value = 3294
$ws.Range("F7").Value = 1129.663694596235
$ws.Range("G7").Value = 20.41398132355528
